$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the "Related Requirement" values (M2:M14) that all repeated "Mina";
# the header in M1 is left in place.
$ws.Range("M2:M14").Clear()

# Add the missing "Car_SRS_35" related-requirement value for row 9, matching the
# formatting already used for the other "Related Requirement" cells in column B.
$ws.Range("B4").Copy()
$ws.Range("B9").PasteSpecial(-4122)
$ws.Range("B9").Value = "Car_SRS_35"

# Update the selection/active cell to match the reviewed state
$ws.Range("B2").Select()
